$wb = $excel.ActiveWorkbook

$wsLibrary  = $wb.Worksheets.Item("Library")
$wsFormula  = $wb.Worksheets.Item("Library_Formula")

# --- Library sheet: drop the trailing blank row (A1:C2 -> A1:C1) ---
$wsLibrary.Rows.Item(2).Delete()

# --- Library_Formula sheet: append the new "getSqlFromExpQuery" formula row ---
$newRow = 5

# Columns A, B, E, F: copy the formatting used by the existing rows above,
# then overwrite with the new row's values.
$wsFormula.Cells.Item(4, 1).Copy()
$wsFormula.Cells.Item($newRow, 1).PasteSpecial(-4122)

$wsFormula.Cells.Item(4, 2).Copy()
$wsFormula.Cells.Item($newRow, 2).PasteSpecial(-4122)

$wsFormula.Cells.Item(4, 5).Copy()
$wsFormula.Cells.Item($newRow, 5).PasteSpecial(-4122)

$wsFormula.Cells.Item(3, 6).Copy()
$wsFormula.Cells.Item($newRow, 6).PasteSpecial(-4122)

$wsFormula.Cells.Item($newRow, 1).Value = "CREATE/MODIFY"
$wsFormula.Cells.Item($newRow, 2).Value = "LIB_VMCommons"
$wsFormula.Cells.Item($newRow, 3).Value = "getSqlFromExpQuery"
$wsFormula.Cells.Item($newRow, 5).Value = "String"
$wsFormula.Cells.Item($newRow, 6).Value = "String, String, String"

# Move the Library_Formula selection down one row (A6 -> A7) ...
$wsFormula.Range("A7").Select()

# ... and make "Library" the active/selected sheet again.
$wsLibrary.Activate()
